$wb = $excel.ActiveWorkbook

# All sheets in this workbook share the label cells in A3 ("premature 2.5")
# and A4 ("premature 3.0"). Update the wording on every sheet so the
# underlying shared-string text is rewritten in place, and make A4 use the
# same number-format style as A3 (style index 4 / "0.0").
$sheetNames = @(
    "Portex uncuffed",
    "Portex ivory",
    "Rusch uncuffed",
    "Sheridan uncuffed",
    "Malinckrodt",
    "Malinckrodt contour"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A3").Value = "premature (expected ID 2.5)"
    $ws.Range("A4").Value = "premature (expected ID 3.0)"

    # Give A4 the same number format as A3 so it picks up style index 4.
    $ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat

    # Select A3:A4 on every sheet (matches the saved selection state).
    $ws.Range("A3:A4").Select()
}

# Make "Portex uncuffed" (the first sheet) the active tab again.
$wb.Worksheets.Item("Portex uncuffed").Activate()
